$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 712816
$ws.Range("C3").Value = 2161000
$ws.Range("C4").Value = 8982000
$ws.Range("C5").Value = 566300000

$ws.Range("J5").Select()
